$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.69"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.297"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05801"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.491"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.336"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8089"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8719"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07282"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03067"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03056"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09317"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.854"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001539"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006012"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006045"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001268"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004589"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.576"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3207"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002345"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03788"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006381"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002700"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006887"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005482"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5502"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.006837"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
